# Auto-generated edit script: updates cryptos list values per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.945.45"
$ws.Range("E2").Value = "  -3.17%  "
$ws.Range("D3").Value = "2.822.21"
$ws.Range("E3").Value = "  -3.75%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Formula = '="502.88"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  -5.55%  "
$ws.Range("D6").Formula = '="134.84"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  -7.93%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Formula = '="0.530"'
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = "  -5.49%  "
$ws.Range("D9").Value = "2.817.30"
$ws.Range("E9").Value = "  -3.36%  "
$ws.Range("D10").Formula = '="5.92"'
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = "  -0.50%  "
$ws.Range("E11").Value = "  -6.76%  "
$ws.Range("D12").Formula = '="0.346"'
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = "  -3.15%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Formula = '="0.127"'
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = "  +1.38%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "3.311.47"
$ws.Range("E14").Value = "  -3.58%  "
$ws.Range("D15").Value = "59.116.61"
$ws.Range("E15").Value = "  -2.98%  "
$ws.Range("D16").Formula = '="21.46"'
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value = "  -7.29%  "
$ws.Range("D17").Value = "2.823.90"
$ws.Range("E17").Value = "  -3.56%  "
$ws.Range("E18").Value = "  -5.98%  "
$ws.Range("E19").Value = "  -6.44%  "
$ws.Range("E20").Value = "  -6.47%  "
$ws.Range("D21").Formula = '="346.82"'
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  -5.22%  "
$ws.Range("D22").Formula = '="6.25"'
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = "  -4.41%  "
$ws.Range("D23").Formula = '="1.00"'
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("D24").Formula = '="5.63"'
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = "  -0.68%  "
$ws.Range("D25").Formula = '="62.91"'
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = "  -2.44%  "
$ws.Range("D26").Formula = '="0.424"'
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = "  -6.94%  "
$ws.Range("E27").Value = "  -6.72%  "
$ws.Range("D28").Formula = '="1.01"'
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = "  +0.92%  "
$ws.Range("D29").Formula = '="7.28"'
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = "  -7.03%  "
$ws.Range("D30").Value = "0.0₃0794"
$ws.Range("E30").Value = "  -10.58%  "
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("E32").Value = "  -5.21%  "
$ws.Range("D33").Formula = '="18.89"'
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = "  -5.03%  "
$ws.Range("D34").Formula = '="150.12"'
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = "  -4.87%  "
$ws.Range("E35").Value = "  -6.99%  "
$ws.Range("D36").Formula = '="5.29"'
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("D37").Formula = '="0.893"'
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = "  -12.44%  "
$ws.Range("E38").Value = "  -8.87%  "
$ws.Range("D39").Formula = '="36.63"'
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = "  -3.24%  "
$ws.Range("D40").Value = "2.213.14"
$ws.Range("E40").Value = "  -6.55%  "
$ws.Range("E41").Value = "  -3.20%  "
$ws.Range("D42").Formula = '="3.50"'
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = "  -6.70%  "
$ws.Range("E43").Value = "  +0.46%  "
$ws.Range("D44").Formula = '="0.0555"'
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = "  -3.68%  "
$ws.Range("E45").Value = "  -10.26%  "
$ws.Range("D46").Formula = '="19.08"'
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = "  -9.61%  "
$ws.Range("D47").Formula = '="10.34"'
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = "  -0.63%  "
$ws.Range("E48").Value = "  -5.08%  "
$ws.Range("D49").Formula = '="0.0887"'
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = "  -4.96%  "
$ws.Range("D50").Formula = '="4.52"'
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = "  -8.45%  "
$ws.Range("E51").Value = "  -7.51%  "
